$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.882332563400269
$ws.Range("B1").Value = 2.00743842124939
$ws.Range("C1").Value = 2.066988706588745
$ws.Range("D1").Value = 2.644529819488525
$ws.Range("E1").Value = 3.811570644378662
